$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 606.5333000000001
$ws.Range("I33").Value = 473.72726
$ws.Range("J33").Value = 971.75
$ws.Range("K33").Value = 473.72726
$ws.Range("L33").Value = 971.75
$ws.Range("M33").Value = -244.72726
$ws.Range("N33").Value = -1429.75
$ws.Range("H100").Value = 2910.7144
$ws.Range("I100").Value = 2347
$ws.Range("K100").Value = 2347
$ws.Range("M100").Value = -1806
$ws.Range("H134").Value = 25656.924
$ws.Range("J134").Value = 25656.924
$ws.Range("L134").Value = 25656.924
$ws.Range("N134").Value = -35796.924
$ws.Range("H135").Value = 609.7037
$ws.Range("I135").Value = 360.09525
$ws.Range("J135").Value = 1483.3334
$ws.Range("K135").Value = 3240.85725
$ws.Range("L135").Value = 13350.0006
$ws.Range("M135").Value = -705.85725
$ws.Range("N135").Value = -18420.0006
$ws.Range("H136").Value = 30000
$ws.Range("J136").Value = 30000
$ws.Range("L136").Value = 30000
$ws.Range("N136").Value = -40200
$ws.Range("H139").Value = 30000
$ws.Range("J139").Value = 30000
$ws.Range("L139").Value = 30000
$ws.Range("N139").Value = -40280
$ws.Range("H140").Value = 25971.428
$ws.Range("J140").Value = 25971.428
$ws.Range("L140").Value = 25971.428
$ws.Range("N140").Value = -36331.428
$ws.Range("H141").Value = 376197.28
$ws.Range("J141").Value = 676226
$ws.Range("L141").Value = 2028678
$ws.Range("N141").Value = -2039038

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 19000
$ws.Range("J9").Value = 19000
$ws.Range("L9").Value = 19000
$ws.Range("N9").Value = -19340
$ws.Range("H20").Value = 19000
$ws.Range("J20").Value = 19000
$ws.Range("L20").Value = 19000
$ws.Range("N20").Value = -19540
$ws.Range("H32").Value = 2776.22
$ws.Range("I32").Value = 2776.22
$ws.Range("K32").Value = 2776.22
$ws.Range("M32").Value = -2489.22
$ws.Range("H37").Value = 17110.4
$ws.Range("J37").Value = 17110.4
$ws.Range("L37").Value = 17110.4
$ws.Range("N37").Value = -17656.4
$ws.Range("H110").Value = 1652.4642
$ws.Range("I110").Value = 640.82355
$ws.Range("J110").Value = 3215.9092
$ws.Range("K110").Value = 640.82355
$ws.Range("L110").Value = 3215.9092
$ws.Range("M110").Value = 1404.17645
$ws.Range("N110").Value = -7305.9092
$ws.Range("H122").Value = 2287
$ws.Range("I122").Value = 1716.8276
$ws.Range("J122").Value = 5594
$ws.Range("K122").Value = 5150.4828
$ws.Range("L122").Value = 16782
$ws.Range("M122").Value = -2700.4828
$ws.Range("N122").Value = -21682
$ws.Range("H132").Value = 2129.547
$ws.Range("I132").Value = 1587.0488
$ws.Range("K132").Value = 4761.1464
$ws.Range("M132").Value = -2231.1464

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 5383
$ws.Range("I99").Value = 5074.5
$ws.Range("K99").Value = 5074.5
$ws.Range("M99").Value = -3576.5
$ws.Range("H107").Value = 2530.077
$ws.Range("I107").Value = 1899.1818
$ws.Range("J107").Value = 6000
$ws.Range("K107").Value = 1899.1818
$ws.Range("L107").Value = 6000
$ws.Range("M107").Value = 20.81819999999993
$ws.Range("N107").Value = -9840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1432.2
$ws.Range("I16").Value = 767.9231
$ws.Range("J16").Value = 5750
$ws.Range("K16").Value = 767.9231
$ws.Range("L16").Value = 5750
$ws.Range("M16").Value = -480.9231
$ws.Range("N16").Value = -6324
$ws.Range("H22").Value = 1067.2
$ws.Range("I22").Value = 295.1111
$ws.Range("J22").Value = 2225.3333
$ws.Range("K22").Value = 295.1111
$ws.Range("L22").Value = 2225.3333
$ws.Range("M22").Value = 54.88889999999998
$ws.Range("N22").Value = -2925.3333
$ws.Range("H53").Value = 27942
$ws.Range("J53").Value = 27942
$ws.Range("L53").Value = 27942
$ws.Range("N53").Value = -29156
$ws.Range("H99").Value = 3586.1667
$ws.Range("I99").Value = 2194.818
$ws.Range("J99").Value = 5772.5713
$ws.Range("K99").Value = 2194.818
$ws.Range("L99").Value = 5772.5713
$ws.Range("M99").Value = -696.8180000000002
$ws.Range("N99").Value = -8768.5713
$ws.Range("H107").Value = 1334
$ws.Range("I107").Value = 1388.5834
$ws.Range("J107").Value = 1283.6154
$ws.Range("K107").Value = 1388.5834
$ws.Range("L107").Value = 1283.6154
$ws.Range("M107").Value = 531.4166
$ws.Range("N107").Value = -5123.6154
$ws.Range("H108").Value = 24614
$ws.Range("I108").Value = 12000
$ws.Range("J108").Value = 27136.8
$ws.Range("K108").Value = 12000
$ws.Range("L108").Value = 27136.8
$ws.Range("M108").Value = -8160
$ws.Range("N108").Value = -34816.8
$ws.Range("H110").Value = 30000
$ws.Range("J110").Value = 30000
$ws.Range("L110").Value = 30000
$ws.Range("N110").Value = -38180
$ws.Range("H111").Value = 30466.666
$ws.Range("J111").Value = 30466.666
$ws.Range("L111").Value = 30466.666
$ws.Range("N111").Value = -38646.666
$ws.Range("H113").Value = 1432.2
$ws.Range("I113").Value = 767.9231
$ws.Range("J113").Value = 5750
$ws.Range("K113").Value = 767.9231
$ws.Range("L113").Value = 5750
$ws.Range("M113").Value = 1402.0769
$ws.Range("N113").Value = -10090
$ws.Range("H126").Value = 3586.1667
$ws.Range("I126").Value = 2194.818
$ws.Range("J126").Value = 5772.5713
$ws.Range("K126").Value = 6584.454000000001
$ws.Range("L126").Value = 17317.7139
$ws.Range("M126").Value = -4114.454000000001
$ws.Range("N126").Value = -22257.7139
$ws.Range("H134").Value = 9806024
$ws.Range("I134").Value = 12501538
$ws.Range("J134").Value = 4154.4546
$ws.Range("K134").Value = 37504614
$ws.Range("L134").Value = 12463.3638
$ws.Range("M134").Value = -37502079
$ws.Range("N134").Value = -17533.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 678.8723
$ws.Range("I5").Value = 414.34885
$ws.Range("J5").Value = 3522.5
$ws.Range("K5").Value = 1243.04655
$ws.Range("L5").Value = 10567.5
$ws.Range("M5").Value = -1131.04655
$ws.Range("N5").Value = -10791.5
$ws.Range("H34").Value = 6639.778
$ws.Range("I34").Value = 136.66667
$ws.Range("J34").Value = 9891.333000000001
$ws.Range("K34").Value = 410.00001
$ws.Range("L34").Value = 29673.999
$ws.Range("M34").Value = -326.00001
$ws.Range("N34").Value = -29841.999
$ws.Range("H38").Value = 411.32
$ws.Range("I38").Value = 80.63636
$ws.Range("J38").Value = 671.1429000000001
$ws.Range("K38").Value = 241.90908
$ws.Range("L38").Value = 2013.4287
$ws.Range("M38").Value = 105.09092
$ws.Range("N38").Value = -2707.4287
$ws.Range("H39").Value = 2726
$ws.Range("J39").Value = 2726
$ws.Range("L39").Value = 8178
$ws.Range("N39").Value = -8766
$ws.Range("H55").Value = 2415.8333
$ws.Range("J55").Value = 3087.3076
$ws.Range("L55").Value = 9261.9228
$ws.Range("N55").Value = -9615.9228
$ws.Range("H98").Value = 264.625
$ws.Range("J98").Value = 406.8
$ws.Range("L98").Value = 1220.4
$ws.Range("N98").Value = -4216.4
$ws.Range("H135").Value = 678.8723
$ws.Range("I135").Value = 414.34885
$ws.Range("J135").Value = 3522.5
$ws.Range("K135").Value = 3729.13965
$ws.Range("L135").Value = 31702.5
$ws.Range("M135").Value = -1194.13965
$ws.Range("N135").Value = -36772.5
$ws.Range("H139").Value = 8335660.5
$ws.Range("I139").Value = 10418577
$ws.Range("K139").Value = 31255731
$ws.Range("M139").Value = -31250591

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1118.8
$ws.Range("I107").Value = 458.2
$ws.Range("K107").Value = 458.2
$ws.Range("M107").Value = 1461.8
$ws.Range("H122").Value = 4007.6182
$ws.Range("I122").Value = 2971.6487
$ws.Range("J122").Value = 6137.1113
$ws.Range("K122").Value = 8914.946100000001
$ws.Range("L122").Value = 18411.3339
$ws.Range("M122").Value = -6464.946100000001
$ws.Range("N122").Value = -23311.3339
$ws.Range("H132").Value = 3033.5925
$ws.Range("I132").Value = 2719.6216
$ws.Range("J132").Value = 3716.9412
$ws.Range("K132").Value = 8158.864799999999
$ws.Range("L132").Value = 11150.8236
$ws.Range("M132").Value = -5628.864799999999
$ws.Range("N132").Value = -16210.8236

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2023.9678
$ws.Range("I132").Value = 1346
$ws.Range("J132").Value = 3347.6191
$ws.Range("K132").Value = 4038
$ws.Range("L132").Value = 10042.8573
$ws.Range("M132").Value = -1508
$ws.Range("N132").Value = -15102.8573

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1171.091
$ws.Range("I107").Value = 284.875
$ws.Range("J107").Value = 3534.3333
$ws.Range("K107").Value = 854.625
$ws.Range("L107").Value = 10602.9999
$ws.Range("M107").Value = 1065.375
$ws.Range("N107").Value = -14442.9999
